$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness column (C) for rows 2 through 32 from 7310 to 7293
$ws.Range("C2:C32").Value = 7293
